$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65:141 down to 66:142.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly price record
# (same market/product metadata as the former row 65, new date + price info).
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 45036
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100104
$ws.Range("H65").Value = "Frutos de pepita"
$ws.Range("I65").Value = 100104003
$ws.Range("J65").Value = "Membrillo"
$ws.Range("K65").Value = "Champion"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 14000
$ws.Range("O65").Value = 15000
$ws.Range("P65").Value = 14500
$ws.Range("Q65").Value = "$/caja 18 kilos empedrada"
$ws.Range("R65").Value = "Región de O'Higgins"
$ws.Range("S65").Value = 806
$ws.Range("T65").Value = 18
